# Applies:
#  1. Bold the "Development log" heading paragraph.
#  2. "I did a bit of research" -> "I did research"
#  3. Replace ", bridge and " with " – bridge - " in the "Finished light box..." paragraph.

$d = $word.ActiveDocument

# --- 1. Bold "Development log" paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n") -eq "Development log") {
        $p.Range.Bold = 1
        break
    }
}

# --- 2. "I did a bit of research" -> "I did research" ---
$d.Content.Find.Execute("I did a bit of research", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "I did research", 2)

# --- 3. Rework the "Finished light box..." paragraph into several runs,
#        swapping the comma/"and" joiners for an en-dash and hyphen, and
#        merging the "(t" / "he vista effect)" split around the bookmark
#        into clean " (the vista effect)" text. ---
$enDash = [char]0x2013
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Finished light box*") {
        $start = $p.Range.Start
        $end = $p.Range.End - 1   # exclude the paragraph mark
        $r = $d.Range($start, $end)

        $newXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Finished light box of the concept Idea. Not going with a U shape for the dungeon, just a step</w:t></w:r><w:r><w:t xml:space="preserve"> $enDash</w:t></w:r><w:r><w:t xml:space="preserve"> bridge</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">- </w:t></w:r><w:r><w:t>step to final room. I revisited one of our class time exercises and tried to bring a sense of wonder to the final room while eluding to it the entire time</w:t></w:r><w:r><w:t xml:space="preserve"> (the vista effect)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
        $r.InsertXML($newXml)
        break
    }
}
